$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need the
# Text number format applied first, otherwise Excel auto-converts them
# to actual numbers (losing the "37.38"-style text formatting used by
# this price column).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D12", "D13", "D15", "D16", "D19", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D41", "D42", "D45", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume values scraped for this run.
$ws.Range("D2").Value = "37.334.25"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.034.08"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "248.57"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D7").Value = "60.79"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("D9").Value = "0.396"
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("D10").Value = "0.0812"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "15.34"
$ws.Range("E12").Value = "  +8.49%  "
$ws.Range("D13").Value = "0.863"
$ws.Range("E13").Value = "  +4.80%  "
$ws.Range("D14").Value = "2.335.32"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("D15").Value = "22.41"
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +6.17%  "
$ws.Range("D17").Value = "2.034.85"
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("D18").Value = "37.291.37"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "70.82"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").Value = "0.0₃0869"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").Value = "5.27"
$ws.Range("E21").Value = "  +4.58%  "
$ws.Range("D22").Value = "231.40"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  +6.83%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "163.89"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("D29").Value = "19.87"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("D32").Value = "4.85"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("D33").Value = "0.0674"
$ws.Range("E33").Value = "  +10.67%  "
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +11.81%  "
$ws.Range("D36").Value = "3.67"
$ws.Range("E36").Value = "  +6.33%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("D41").Value = "0.0983"
$ws.Range("E41").Value = "  +3.31%  "
$ws.Range("D42").Value = "17.21"
$ws.Range("E42").Value = "  +10.11%  "
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "92.85"
$ws.Range("E45").Value = "  +5.28%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("D47").Value = "1.388.79"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").Value = "7.50"
$ws.Range("E48").Value = "  +5.98%  "
$ws.Range("D49").Value = "2.18"
$ws.Range("E49").Value = "  +20.83%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "46.47"
$ws.Range("E51").Value = "  +2.94%  "
